$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 197. This pushes the existing rows 197-198
# (the week-of-2021-12-22 data) down to 199-200, preserving their values and
# formatting (incl. the date cell's custom number format).
$ws.Rows("197:198").Insert()

# Fill in the new week's data (2022-02-03 = serial 44595) into the freshly
# inserted rows 197-198, using the same layout as the rows that follow.

# Row 197 - "Primera" quality
$ws.Cells.Item(197, 1).Value = 11
$ws.Cells.Item(197, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(197, 3).Value = "Bíobío"
$ws.Cells.Item(197, 4).Value = 44595
$ws.Cells.Item(197, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(197, 5).Value = 8
$ws.Cells.Item(197, 6).Value = 100112017
$ws.Cells.Item(197, 7).Value = "Apio"
$ws.Cells.Item(197, 8).Value = "Americana (o)"
$ws.Cells.Item(197, 9).Value = "Primera"
$ws.Cells.Item(197, 10).Value = 100
$ws.Cells.Item(197, 11).Value = 6000
$ws.Cells.Item(197, 12).Value = 7000
$ws.Cells.Item(197, 13).Value = 6500
$ws.Cells.Item(197, 14).Value = "$/docena de matas"
$ws.Cells.Item(197, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(197, 16).Value = 1083
$ws.Cells.Item(197, 17).Value = 6
$ws.Cells.Item(197, 18).Value = "Hortaliza"

# Row 198 - "Segunda" quality
$ws.Cells.Item(198, 1).Value = 11
$ws.Cells.Item(198, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(198, 3).Value = "Bíobío"
$ws.Cells.Item(198, 4).Value = 44595
$ws.Cells.Item(198, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(198, 5).Value = 8
$ws.Cells.Item(198, 6).Value = 100112017
$ws.Cells.Item(198, 7).Value = "Apio"
$ws.Cells.Item(198, 8).Value = "Americana (o)"
$ws.Cells.Item(198, 9).Value = "Segunda"
$ws.Cells.Item(198, 10).Value = 50
$ws.Cells.Item(198, 11).Value = 5000
$ws.Cells.Item(198, 12).Value = 5000
$ws.Cells.Item(198, 13).Value = 5000
$ws.Cells.Item(198, 14).Value = "$/docena de matas"
$ws.Cells.Item(198, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(198, 16).Value = 833
$ws.Cells.Item(198, 17).Value = 6
$ws.Cells.Item(198, 18).Value = "Hortaliza"
